$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $parts = $val.Split(",")
    $count = $parts.Length

    if ($count -gt 1) {
        for ($i = 0; $i -lt $count; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        $reversed = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }

        $newVal = $reversed -join ", "
        $cell.Value2 = $newVal
    }
}
